$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = 1050000
$ws.Range("F9").Value = 1050000
$ws.Range("D12").Value = 1403940
$ws.Range("F12").Value = 1403940
$ws.Range("D16").Value = 1709940
$ws.Range("F16").Value = 1709940
$ws.Range("D20").Value = 2507940
$ws.Range("F20").Value = 2507940
$ws.Range("D21").Value = 1397940
$ws.Range("F21").Value = 1397940
$ws.Range("D23").Value = 1673940
$ws.Range("F23").Value = 1673940
$ws.Range("D24").Value = 2039940
$ws.Range("F24").Value = 2039940
$ws.Range("D26").Value = 1389000
$ws.Range("F26").Value = 1389000
$ws.Range("D27").Value = 1655940
$ws.Range("F27").Value = 1655940
$ws.Range("D28").Value = 1722000
$ws.Range("F28").Value = 1722000
$ws.Range("D30").Value = 1398000
$ws.Range("F30").Value = 1398000
$ws.Range("D34").Value = 802900
$ws.Range("F34").Value = 802900
$ws.Range("D40").Value = 2850000
$ws.Range("F40").Value = 2850000
$ws.Range("D65").Value = 1976394
$ws.Range("F65").Value = 1976394
$ws.Range("D67").Value = 1619994
$ws.Range("F67").Value = 1619994
$ws.Range("D68").Value = 963300
$ws.Range("F68").Value = 963300
$ws.Range("D69").Value = 1583994
$ws.Range("F69").Value = 1583994
$ws.Range("D70").Value = 1267794
$ws.Range("F70").Value = 1267794
$ws.Range("D72").Value = 2596194
$ws.Range("F72").Value = 2596194
$ws.Range("D73").Value = 2117394
$ws.Range("F73").Value = 2117394
$ws.Range("D79").Value = 1687000
$ws.Range("F79").Value = 1687000
$ws.Range("D80").Value = 2009400
$ws.Range("F80").Value = 2009400
$ws.Range("D82").Value = 1619400
$ws.Range("F82").Value = 1619400
$ws.Range("D88").Value = 1199400
$ws.Range("F88").Value = 1199400
